# feat: add 2022-Q3 data
#
# Current workbook layout:
#   Sheet 1: "总计"     (summary)
#   Sheet 2: "2022-Q2"  (quarterly fund-holding detail)
#
# Target workbook layout:
#   Sheet 1: "总计"     (summary, gains a new row for 2022-Q3, old row becomes 2022-Q2)
#   Sheet 2: "2022-Q3"  (NEW quarterly fund-holding detail)
#   Sheet 3: "2022-Q2"  (the OLD quarterly fund-holding detail, moved down)

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (so that numeric-looking strings such as "168105" or "0.51" are not
# silently converted into numbers), then clear the number-format override
# so no stray style survives on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet so its data/format is
# preserved; the duplicate is placed right after it and becomes the new
# active sheet, matching how Excel itself behaves on Sheet.Copy.
# ---------------------------------------------------------------------
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Step 2: rename sheets.  The original sheet becomes "2022-Q3" (it keeps
# sheetId 2 / rId2); the duplicate becomes the new "2022-Q2" (sheetId 3 /
# rId3).  Original must be renamed first so the duplicate can take the
# now-freed "2022-Q2" name.
# ---------------------------------------------------------------------
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

$q3 = $q2

# ---------------------------------------------------------------------
# Step 3: replace the contents of the (renamed) "2022-Q3" sheet with the
# brand-new 2022-Q3 fund-holding data.
# ---------------------------------------------------------------------
$q3.Cells.Clear()

# Reset page margins to the workbook's "fresh sheet" defaults (0.75in /
# 0.75in / 1in / 1in / 0.5in / 0.5in), matching a newly created sheet
# instead of the ones inherited from the duplicated "2022-Q2" sheet.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row - copy the "总计" header style (bold / centered / thin
# border, style index 2) across B1:H1, then fill in the Q3 labels.
$total.Range("B1").Copy($q3.Range("B1:H1"))
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2 - fund 168105
$total.Range("A2").Copy($q3.Range("A2"))
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "168105"
$q3.Range("C2").Value = "九泰泰富灵活配置混合（LOF）A"
Set-TextValue $q3.Range("D2") "0.51"
Set-TextValue $q3.Range("E2") "88.30"
Set-TextValue $q3.Range("F2") "3.42"
Set-TextValue $q3.Range("G2") "0.0174"
$q3.Range("H2").Value = 4

# Row 3 - fund 015688
$total.Range("A2").Copy($q3.Range("A3"))
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "015688"
$q3.Range("C3").Value = "九泰泰富灵活配置混合（LOF）C"
Set-TextValue $q3.Range("D3") "0.00"
Set-TextValue $q3.Range("E3") "88.30"
Set-TextValue $q3.Range("F3") "3.42"
$q3.Range("G3").Value = 0
$q3.Range("H3").Value = 4

# ---------------------------------------------------------------------
# Step 4: update the "总计" sheet - row 2 now describes 2022-Q3, and a
# new row 3 is added describing 2022-Q2 (what row 2 used to hold).
# ---------------------------------------------------------------------
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.02
